$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the 8 feature activation bits (row 4, columns B..I)
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 1

# Update the equivalent combined features code (O4)
$ws.Range("O4").Value = 253

# Row 12 (CHECK_EMPTY_FL) J/M cells become the new literal label "APP_AUTORUN"
# instead of the shared formula referencing $S12
$ws.Range("J12").Value = "APP_AUTORUN"
$ws.Range("M12").Value = "APP_AUTORUN"
